$wb = $excel.ActiveWorkbook

# --- 1. Update selection on the first sheet ("Experiments and Notes"):
#     activeCell moves from P31 to N25 (frozen-pane "topRight" selection).
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$ws1.Range("N25").Select()

# --- 2. Add the new "Conversion Rates" worksheet after "Lab Modules" (becomes sheet 4).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Conversion Rates"

# --- 3. Column widths for the new sheet.
$ws.Columns.Item(1).ColumnWidth = 12.96
$ws.Columns.Item(2).ColumnWidth = 15.84
$ws.Columns.Item(3).ColumnWidth = 15.28
$ws.Columns.Item(4).ColumnWidth = 24.45
$ws.Columns.Item(5).ColumnWidth = 13.65

# --- 4. Header row (new shared strings get appended automatically).
$ws.Cells.Item(1,1).Value = "Unit per Hour"
$ws.Cells.Item(1,2).Value = "per second"
$ws.Cells.Item(1,3).Value = "per minute"
$ws.Cells.Item(1,4).Value = "per hour"
$ws.Cells.Item(1,5).Value = "per Kerbin day"

# --- 5. Data rows 2-50 (values + formulas).
$ws.Cells.Item(2,1).Value = 0.1
$ws.Cells.Item(3,1).Value = 0.2
$ws.Cells.Item(4,1).Value = 0.3
$ws.Cells.Item(5,1).Value = 0.4
$ws.Cells.Item(6,1).Value = 0.5
$ws.Cells.Item(7,1).Value = 0.6
$ws.Cells.Item(8,1).Value = 0.7
$ws.Cells.Item(9,1).Value = 0.8
$ws.Cells.Item(10,1).Value = 0.9
$ws.Cells.Item(11,1).Value = 1
$ws.Cells.Item(11,2).Value = 0.00027777777
$ws.Cells.Item(11,3).Value = 0.0166666662
$ws.Cells.Item(11,4).Value = 0.999999972
$ws.Cells.Item(11,5).Value = 3.999999888
$ws.Cells.Item(12,1).Value = 2
$ws.Cells.Item(12,2).Formula = "=B11*A12"
$ws.Cells.Item(12,3).Formula = "=C11*A12"
$ws.Cells.Item(12,4).Formula = "=D11*A12"
$ws.Cells.Item(12,5).Formula = "=E11*A12"
$ws.Cells.Item(13,1).Value = 3
$ws.Cells.Item(13,2).Formula = "=B12*A13"
$ws.Cells.Item(13,3).Formula = "=C12*A13"
$ws.Cells.Item(13,4).Formula = "=D12*A13"
$ws.Cells.Item(13,5).Formula = "=E12*A13"
$ws.Cells.Item(14,1).Value = 4
$ws.Cells.Item(14,2).Formula = "=B13*A14"
$ws.Cells.Item(14,3).Formula = "=C13*A14"
$ws.Cells.Item(14,4).Formula = "=D13*A14"
$ws.Cells.Item(14,5).Formula = "=E13*A14"
$ws.Cells.Item(15,1).Value = 5
$ws.Cells.Item(15,2).Formula = "=B14*A15"
$ws.Cells.Item(15,3).Formula = "=C14*A15"
$ws.Cells.Item(15,4).Formula = "=D14*A15"
$ws.Cells.Item(15,5).Formula = "=E14*A15"
$ws.Cells.Item(16,1).Value = 6
$ws.Cells.Item(16,2).Formula = "=B15*A16"
$ws.Cells.Item(16,3).Formula = "=C15*A16"
$ws.Cells.Item(16,4).Formula = "=D15*A16"
$ws.Cells.Item(16,5).Formula = "=E15*A16"
$ws.Cells.Item(17,1).Value = 7
$ws.Cells.Item(17,2).Formula = "=B16*A17"
$ws.Cells.Item(17,3).Formula = "=C16*A17"
$ws.Cells.Item(17,4).Formula = "=D16*A17"
$ws.Cells.Item(17,5).Formula = "=E16*A17"
$ws.Cells.Item(18,1).Value = 8
$ws.Cells.Item(18,2).Formula = "=B17*A18"
$ws.Cells.Item(18,3).Formula = "=C17*A18"
$ws.Cells.Item(18,4).Formula = "=D17*A18"
$ws.Cells.Item(18,5).Formula = "=E17*A18"
$ws.Cells.Item(19,1).Value = 9
$ws.Cells.Item(19,2).Formula = "=B18*A19"
$ws.Cells.Item(19,3).Formula = "=C18*A19"
$ws.Cells.Item(19,4).Formula = "=D18*A19"
$ws.Cells.Item(19,5).Formula = "=E18*A19"
$ws.Cells.Item(20,1).Value = 10
$ws.Cells.Item(20,2).Formula = "=B19*A20"
$ws.Cells.Item(20,3).Formula = "=C19*A20"
$ws.Cells.Item(20,4).Formula = "=D19*A20"
$ws.Cells.Item(20,5).Formula = "=E19*A20"
$ws.Cells.Item(21,1).Value = 11
$ws.Cells.Item(21,2).Formula = "=B20*A21"
$ws.Cells.Item(21,3).Formula = "=C20*A21"
$ws.Cells.Item(21,4).Formula = "=D20*A21"
$ws.Cells.Item(21,5).Formula = "=E20*A21"
$ws.Cells.Item(22,1).Value = 12
$ws.Cells.Item(22,2).Formula = "=B21*A22"
$ws.Cells.Item(22,3).Formula = "=C21*A22"
$ws.Cells.Item(22,4).Formula = "=D21*A22"
$ws.Cells.Item(22,5).Formula = "=E21*A22"
$ws.Cells.Item(23,1).Value = 13
$ws.Cells.Item(23,2).Formula = "=B22*A23"
$ws.Cells.Item(23,3).Formula = "=C22*A23"
$ws.Cells.Item(23,4).Formula = "=D22*A23"
$ws.Cells.Item(23,5).Formula = "=E22*A23"
$ws.Cells.Item(24,1).Value = 14
$ws.Cells.Item(24,2).Formula = "=B23*A24"
$ws.Cells.Item(24,3).Formula = "=C23*A24"
$ws.Cells.Item(24,4).Formula = "=D23*A24"
$ws.Cells.Item(24,5).Formula = "=E23*A24"
$ws.Cells.Item(25,1).Value = 15
$ws.Cells.Item(25,2).Formula = "=B24*A25"
$ws.Cells.Item(25,3).Formula = "=C24*A25"
$ws.Cells.Item(25,4).Formula = "=D24*A25"
$ws.Cells.Item(25,5).Formula = "=E24*A25"
$ws.Cells.Item(26,1).Value = 16
$ws.Cells.Item(26,2).Formula = "=B25*A26"
$ws.Cells.Item(26,3).Formula = "=C25*A26"
$ws.Cells.Item(26,4).Formula = "=D25*A26"
$ws.Cells.Item(26,5).Formula = "=E25*A26"
$ws.Cells.Item(27,1).Value = 17
$ws.Cells.Item(27,2).Formula = "=B26*A27"
$ws.Cells.Item(27,3).Formula = "=C26*A27"
$ws.Cells.Item(27,4).Formula = "=D26*A27"
$ws.Cells.Item(27,5).Formula = "=E26*A27"
$ws.Cells.Item(28,1).Value = 18
$ws.Cells.Item(28,2).Formula = "=B27*A28"
$ws.Cells.Item(28,3).Formula = "=C27*A28"
$ws.Cells.Item(28,4).Formula = "=D27*A28"
$ws.Cells.Item(28,5).Formula = "=E27*A28"
$ws.Cells.Item(29,1).Value = 19
$ws.Cells.Item(29,2).Formula = "=B28*A29"
$ws.Cells.Item(29,3).Formula = "=C28*A29"
$ws.Cells.Item(29,4).Formula = "=D28*A29"
$ws.Cells.Item(29,5).Formula = "=E28*A29"
$ws.Cells.Item(30,1).Value = 20
$ws.Cells.Item(30,2).Formula = "=B29*A30"
$ws.Cells.Item(30,3).Formula = "=C29*A30"
$ws.Cells.Item(30,4).Formula = "=D29*A30"
$ws.Cells.Item(30,5).Formula = "=E29*A30"
$ws.Cells.Item(31,1).Value = 21
$ws.Cells.Item(31,2).Formula = "=B30*A31"
$ws.Cells.Item(31,3).Formula = "=C30*A31"
$ws.Cells.Item(31,4).Formula = "=D30*A31"
$ws.Cells.Item(31,5).Formula = "=E30*A31"
$ws.Cells.Item(32,1).Value = 22
$ws.Cells.Item(32,2).Formula = "=B31*A32"
$ws.Cells.Item(32,3).Formula = "=C31*A32"
$ws.Cells.Item(32,4).Formula = "=D31*A32"
$ws.Cells.Item(32,5).Formula = "=E31*A32"
$ws.Cells.Item(33,1).Value = 23
$ws.Cells.Item(33,2).Formula = "=B32*A33"
$ws.Cells.Item(33,3).Formula = "=C32*A33"
$ws.Cells.Item(33,4).Formula = "=D32*A33"
$ws.Cells.Item(33,5).Formula = "=E32*A33"
$ws.Cells.Item(34,1).Value = 24
$ws.Cells.Item(34,2).Formula = "=B33*A34"
$ws.Cells.Item(34,3).Formula = "=C33*A34"
$ws.Cells.Item(34,4).Formula = "=D33*A34"
$ws.Cells.Item(34,5).Formula = "=E33*A34"
$ws.Cells.Item(35,1).Value = 25
$ws.Cells.Item(35,2).Formula = "=B34*A35"
$ws.Cells.Item(35,3).Formula = "=C34*A35"
$ws.Cells.Item(35,4).Formula = "=D34*A35"
$ws.Cells.Item(35,5).Formula = "=E34*A35"
$ws.Cells.Item(36,1).Value = 26
$ws.Cells.Item(36,2).Formula = "=B35*A36"
$ws.Cells.Item(36,3).Formula = "=C35*A36"
$ws.Cells.Item(36,4).Formula = "=D35*A36"
$ws.Cells.Item(36,5).Formula = "=E35*A36"
$ws.Cells.Item(37,1).Value = 27
$ws.Cells.Item(37,2).Formula = "=B36*A37"
$ws.Cells.Item(37,3).Formula = "=C36*A37"
$ws.Cells.Item(37,4).Formula = "=D36*A37"
$ws.Cells.Item(37,5).Formula = "=E36*A37"
$ws.Cells.Item(38,1).Value = 28
$ws.Cells.Item(38,2).Formula = "=B37*A38"
$ws.Cells.Item(38,3).Formula = "=C37*A38"
$ws.Cells.Item(38,4).Formula = "=D37*A38"
$ws.Cells.Item(38,5).Formula = "=E37*A38"
$ws.Cells.Item(39,1).Value = 29
$ws.Cells.Item(39,2).Formula = "=B38*A39"
$ws.Cells.Item(39,3).Formula = "=C38*A39"
$ws.Cells.Item(39,4).Formula = "=D38*A39"
$ws.Cells.Item(39,5).Formula = "=E38*A39"
$ws.Cells.Item(40,1).Value = 30
$ws.Cells.Item(40,2).Formula = "=B39*A40"
$ws.Cells.Item(40,3).Formula = "=C39*A40"
$ws.Cells.Item(40,4).Formula = "=D39*A40"
$ws.Cells.Item(40,5).Formula = "=E39*A40"
$ws.Cells.Item(41,1).Value = 31
$ws.Cells.Item(41,2).Formula = "=B40*A41"
$ws.Cells.Item(41,3).Formula = "=C40*A41"
$ws.Cells.Item(41,4).Formula = "=D40*A41"
$ws.Cells.Item(41,5).Formula = "=E40*A41"
$ws.Cells.Item(42,1).Value = 32
$ws.Cells.Item(42,2).Formula = "=B41*A42"
$ws.Cells.Item(42,3).Formula = "=C41*A42"
$ws.Cells.Item(42,4).Formula = "=D41*A42"
$ws.Cells.Item(42,5).Formula = "=E41*A42"
$ws.Cells.Item(43,1).Value = 33
$ws.Cells.Item(43,2).Formula = "=B42*A43"
$ws.Cells.Item(43,3).Formula = "=C42*A43"
$ws.Cells.Item(43,4).Formula = "=D42*A43"
$ws.Cells.Item(43,5).Formula = "=E42*A43"
$ws.Cells.Item(44,1).Value = 34
$ws.Cells.Item(44,2).Formula = "=B43*A44"
$ws.Cells.Item(44,3).Formula = "=C43*A44"
$ws.Cells.Item(44,4).Formula = "=D43*A44"
$ws.Cells.Item(44,5).Formula = "=E43*A44"
$ws.Cells.Item(45,1).Value = 35
$ws.Cells.Item(45,2).Formula = "=B44*A45"
$ws.Cells.Item(45,3).Formula = "=C44*A45"
$ws.Cells.Item(45,4).Formula = "=D44*A45"
$ws.Cells.Item(45,5).Formula = "=E44*A45"
$ws.Cells.Item(46,1).Value = 36
$ws.Cells.Item(46,2).Formula = "=B45*A46"
$ws.Cells.Item(46,3).Formula = "=C45*A46"
$ws.Cells.Item(46,4).Formula = "=D45*A46"
$ws.Cells.Item(46,5).Formula = "=E45*A46"
$ws.Cells.Item(47,1).Value = 37
$ws.Cells.Item(47,2).Formula = "=B46*A47"
$ws.Cells.Item(47,3).Formula = "=C46*A47"
$ws.Cells.Item(47,4).Formula = "=D46*A47"
$ws.Cells.Item(47,5).Formula = "=E46*A47"
$ws.Cells.Item(48,1).Value = 38
$ws.Cells.Item(48,2).Formula = "=B47*A48"
$ws.Cells.Item(48,3).Formula = "=C47*A48"
$ws.Cells.Item(48,4).Formula = "=D47*A48"
$ws.Cells.Item(48,5).Formula = "=E47*A48"
$ws.Cells.Item(49,1).Value = 39
$ws.Cells.Item(49,2).Formula = "=B48*A49"
$ws.Cells.Item(49,3).Formula = "=C48*A49"
$ws.Cells.Item(49,4).Formula = "=D48*A49"
$ws.Cells.Item(49,5).Formula = "=E48*A49"
$ws.Cells.Item(50,1).Value = 40
$ws.Cells.Item(50,2).Formula = "=B49*A50"
$ws.Cells.Item(50,3).Formula = "=C49*A50"
$ws.Cells.Item(50,4).Formula = "=D49*A50"
$ws.Cells.Item(50,5).Formula = "=E49*A50"

# --- 6. Give B11:C11 (the "1 unit per hour" reference values) the Times New
#     Roman 12pt styling used for the callout numbers.
$ref = $ws.Range("B11:C11")
$ref.Font.Name = "Times New Roman"
$ref.Font.Size = 12

# --- 7. Make the new sheet the active tab/selection, matching the saved view.
$ws.Activate()
$ws.Range("H17").Select()
